# "6 hours by turn fix" - MEC-1A schedule sheet
# The class/break schedule is shifted/expanded to fit 6 hours per turn:
# lunch now starts at 12:20 (row 9) instead of 11:30 (old row 8), every
# later row slides down by one, and three new 50-min slots (16:40, 17:30,
# 18:20) are appended at the bottom. A few weekday assignment cells in the
# morning block (rows 3,4,6,7) are also corrected to reflect the new turn
# layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Morning block corrections (7:50 / 8:40 / 9:50 / 10:40 rows) ---
$ws.Range('C3').Value = '-'
$ws.Range('D3').Value = 'Euclides-Gestão'
$ws.Range('E3').Value = '-'

$ws.Range('D4').Value = 'João Rodrigues-Desenho Técnico'
$ws.Range('E4').Value = '-'
$ws.Range('F4').Value = '-'

$ws.Range('C6').Value = 'Andre Lucca-Circuitos Elétricos'
$ws.Range('D6').Value = 'João Rodrigues-Desenho Técnico'
$ws.Range('E6').Value = 'José Ferreira-Tecnologia dos Materiais'
$ws.Range('F6').Value = '-'

$ws.Range('C7').Value = 'José Ferreira-Tecnologia dos Materiais'
$ws.Range('E7').Value = 'José Ferreira-Tecnologia dos Materiais'
$ws.Range('F7').Value = 'Andre Lucca-Circuitos Elétricos'

# --- Old lunch row (11:30) becomes a regular empty class slot ---
$ws.Range('B8').Value = '-'
$ws.Range('C8').Value = '-'
$ws.Range('D8').Value = '-'
$ws.Range('E8').Value = '-'
$ws.Range('F8').Value = '-'

# --- Lunch moves to 12:20 (row 9) ---
$ws.Range('A9').Value = '12:20'
$ws.Range('B9').Value = 'Almoço'
$ws.Range('C9').Value = 'Almoço'
$ws.Range('D9').Value = 'Almoço'
$ws.Range('E9').Value = 'Almoço'
$ws.Range('F9').Value = 'Almoço'

# --- Afternoon time labels shift down one slot (values stay '-') ---
$ws.Range('A10').Value = '13:00'
$ws.Range('A11').Value = '13:50'
$ws.Range('A12').Value = '14:40'
$ws.Range('B12').Value = '-'
$ws.Range('C12').Value = '-'
$ws.Range('D12').Value = '-'
$ws.Range('E12').Value = '-'
$ws.Range('F12').Value = '-'

# --- Afternoon break shifts down to row 13 (15:30) ---
$ws.Range('A13').Value = '15:30'
$ws.Range('B13').Value = 'Intervalo'
$ws.Range('C13').Value = 'Intervalo'
$ws.Range('D13').Value = 'Intervalo'
$ws.Range('E13').Value = 'Intervalo'
$ws.Range('F13').Value = 'Intervalo'

# --- Row 14 keeps '-' values, only its time label shifts ---
$ws.Range('A14').Value = '15:50'

# --- New rows 15-17 appended for the extra turn hours ---
$ws.Range('A15').Value = '16:40'
$ws.Range('B15').Value = '-'
$ws.Range('C15').Value = '-'
$ws.Range('D15').Value = '-'
$ws.Range('E15').Value = '-'
$ws.Range('F15').Value = '-'

$ws.Range('A16').Value = '17:30'
$ws.Range('B16').Value = '-'
$ws.Range('C16').Value = '-'
$ws.Range('D16').Value = '-'
$ws.Range('E16').Value = '-'
$ws.Range('F16').Value = '-'

$ws.Range('A17').Value = '18:20'
$ws.Range('B17').Value = ''
$ws.Range('C17').Value = ''
$ws.Range('D17').Value = ''
$ws.Range('E17').Value = ''
$ws.Range('F17').Value = ''
